$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @{
    2  = "DTaP/"
    3  = "DTaP/"
    4  = "DTaP/"
    5  = "DTaP-Hib "
    6  = "e-IPV"
    7  = "Hepatitis B-Hibi "
    8  = "Hepatitis A Pediatric"
    9  = "Hepatitis A Pediatric"
    10 = "Hepatitis A Adult"
    11 = "Hepatitis A Adult"
    12 = "Hepatitis B Preservative Free Pediatric/Adolescent"
    13 = "Hepatitis B Preservative Free Pediatric/Adolescent"
    14 = "Hepatitis B 2 dose Adolescent (11-15)"
    15 = "Hepatitis B-Adult"
    16 = "Hepatitis B-Adult"
    17 = "Hib"
    18 = "Hib"
    19 = "Hib"
    22 = "MMR/"
    23 = "Measles"
    24 = "Mumps"
    25 = "Pneumococcal 7-valent (Pediatric)"
    26 = "Rubella"
    27 = "Varicella"
}

$colD = @{
    3  = "15 dose vial 10 x 1 dose vial"
    6  = "10 dose vial 10 x 1 dose syringe"
    8  = "5 x 1 dose vial 5 x 1 dose syringes"
    9  = "10 x 1 dose vial 1 dose vial 5 x 1 dose TipLok 25 x 1 dose TipLok"
    10 = "5 x 1 dose vials 5 x 1 dose syringe"
    11 = "1 dose vial 1 dose TIP-LOK 5 x 1 dose TIP-LOK"
    12 = "10 x 1 dose vials 5 x 1 dose 23G TipLok 5 x 1 dose 25G TipLok 251 dose 23G TipLok 251 dose 25G TipLok"
    14 = "10 x 1 dose vials 5 x 1 dose syringes"
    15 = "10 x 1 dose vials 10 x 3 dose vial 5 x 1 dose syringe"
    16 = "1 x 1 dose vial 5 x 1 dose Tiplok 25 x 1 dose Tiplok"
}

foreach ($row in $colA.Keys) {
    $ws.Cells.Item($row, 1).Value = $colA[$row]
}

foreach ($row in $colD.Keys) {
    $ws.Cells.Item($row, 4).Value = $colD[$row]
}
